$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Header date
Replace-Text "2024-03-18 Monday" "2024-03-19 Tuesday"

# Row 1
Replace-Text "34÷8=4, 2" "51÷6=8, 3"
Replace-Text "78÷3=26, 0" "93÷8=11, 5"
Replace-Text "31÷6=5, 1" "39÷2=19, 1"
Replace-Text "76÷4=19, 0" "20÷6=3, 2"
Replace-Text "68÷7=9, 5" "33÷4=8, 1"

# Row 2
Replace-Text "66÷6=11, 0" "41÷5=8, 1"
Replace-Text "12÷8=1, 4" "26÷8=3, 2"
Replace-Text "66÷7=9, 3" "31÷6=5, 1"
Replace-Text "72÷3=24, 0" "42÷7=6, 0"
Replace-Text "39÷4=9, 3" "64÷5=12, 4"

# Row 3
Replace-Text "77÷5=15, 2" "95÷4=23, 3"
Replace-Text "17÷7=2, 3" "99÷2=49, 1"
# Cells 3/4/5 effectively shift values: handle carefully to avoid collision
# (cell4 takes the old cell5 value, cell5 gets a brand-new value)
Replace-Text "28÷4=7, 0" "52÷4=13, 0"
Replace-Text "89÷5=17, 4" "29÷3=9, 2"
Replace-Text "44÷8=5, 4" "28÷4=7, 0"

# Row 4
Replace-Text "60÷2=30, 0" "30÷5=6, 0"
Replace-Text "96÷5=19, 1" "60÷8=7, 4"
Replace-Text "42÷5=8, 2" "39÷6=6, 3"
Replace-Text "82÷6=13, 4" "54÷5=10, 4"
Replace-Text "41÷8=5, 1" "62÷3=20, 2"

# Row 5
Replace-Text "11÷9=1, 2" "68÷9=7, 5"
Replace-Text "12÷3=4, 0" "33÷6=5, 3"
Replace-Text "50÷4=12, 2" "20÷2=10, 0"
Replace-Text "71÷7=10, 1" "55÷9=6, 1"
Replace-Text "90÷2=45, 0" "34÷5=6, 4"
